$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.950639009475708
$ws.Range("B1").Value = 2.244920969009399
$ws.Range("C1").Value = 2.937027215957642
$ws.Range("D1").Value = 5.962118625640869
$ws.Range("E1").Value = 2.961672306060791
